# Release-Notes.xlsx update — folder inventory refresh
# (Folder inventory updated on Mon Jun 16 05:29:54 UTC 2025)

$wb = $excel.ActiveWorkbook

# --- "Folder Inventory" sheet -------------------------------------------
# The newest folder ("Power_Platform_Workshop:Administration_and_Governance")
# moved up to row 2 (with a refreshed "Last Updated" stamp), and the folder
# that used to be in row 2 ("Getting_started_with_Azure_AI_services") is now
# row 3, carrying the timestamp that used to belong to the row-2 entry.
$wsInventory = $wb.Worksheets.Item("Folder Inventory")

$wsInventory.Range("A2").Value = "Power_Platform_Workshop:Administration_and_Governance"
$wsInventory.Range("B2").Value = "Power_Platform_Workshop:Administration_and_Governance"
$wsInventory.Range("C2").Value = "2025-06-16 10:59:35 +0530"

$wsInventory.Range("A3").Value = "Getting_started_with_Azure_AI_services"
$wsInventory.Range("B3").Value = "Getting_started_with_Azure_AI_services"
$wsInventory.Range("C3").Value = "2025-06-16 10:58:16 +0530"

# --- "Metadata" sheet -----------------------------------------------------
$wsMetadata = $wb.Worksheets.Item("Metadata")

# Generated On timestamp refresh
$wsMetadata.Range("B3").Value = "2025-06-16 05:29:54 UTC"

# Workflow Run counter increments 9 -> 10, stored as text (matches the
# existing text-typed "9" cell) rather than being coerced to a number.
$wsMetadata.Range("B5").Value = "'10"

# --- "Summary" sheet --------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

# Most Recent Update mirrors the new row-2 "Last Updated" value.
$wsSummary.Range("B5").Value = "2025-06-16 10:59:35 +0530"
